$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the merged header "Определение интервалов" from V1:W1 to V1:Y1
$ws.Range("V1:W1").UnMerge()
$ws.Range("V1:Y1").Merge()

# New column headers (row 2)
$ws.Range("X2").Value = "Кол-во интервалов"
$ws.Range("Y2").Value = "Ширина"

# New data / formula row (row 3): number of intervals and interval width
$ws.Range("X3").Value = 10
$ws.Range("Y3").Formula = "=(W3-V3)/X3"

# Column width for the new column X
$ws.Columns.Item(24).ColumnWidth = 16.6

# Update active selection to match target state
$ws.Range("Y6").Select()
